# Phone number validation implemented
#
# Registration data: A20's MIS number is normalized from text to a real
# number, and four new registration rows (21-24) are appended - including
# two rows (23-24) whose MIS numbers look like phone numbers but must be
# kept stored as text per the new phone-number validation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: MIS number was stored as text ("612203127"); store it as a
#     real number now that validation normalizes numeric-looking input.
$ws.Range("A20").Value = 612203127

# --- New row 21 (phone number re-used from row 19 but kept as text) ---
$ws.Range("A21").Value = 1234
$ws.Range("B21").Value = "Smruti"
$ws.Range("C21").Value = "a"
$ws.Range("D21").Value = "xyz@gmail.com"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1654654"

# --- New row 22 ---
$ws.Range("A22").Value = "a"
$ws.Range("B22").Value = "a"
$ws.Range("C22").Value = "a"
$ws.Range("D22").Value = "admin@sms.com"
$ws.Range("E22").Value = "a"

# --- New row 23: MIS number & phone look numeric but are kept as text ---
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "112103127"
$ws.Range("B23").Value = "A"
$ws.Range("C23").Value = "B"
$ws.Range("D23").Value = "abc@mail.com"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "987654321"

# --- New row 24: same situation ---
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "612203217"
$ws.Range("B24").Value = "jabsdi"
$ws.Range("C24").Value = "jkhw"
$ws.Range("D24").Value = "abc@mail.com"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "987654321"
